$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update DATE END value (B2) from 14/08/2025 to 25/08/2025
# A2 (DATE START, 20/06/2025) stays the same value.
$ws.Range("B2").Value = "25/08/2025"

# Update the selected cell to B3 (as reflected in the saved view state)
$ws.Range("B3").Select()
